# Applies the "tutorial model in process" edit to simulator_scenarios.xlsx:
#  - inserts a new "period_hours" column (D) after "periods"
#  - inserts a new "young_percentage" column (I) after "initial_infected_percentage"
#  - renames "infection_probability" -> "infection_prob"
#  - appends two new trailing columns: "reinfection_prob" and "vaccinated_infection_prob"
#  - updates column widths / sheet view to the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert "period_hours" as new column D (after "periods", before "agent_num") ---
$ws.Columns.Item(4).Insert()
$ws.Range("D1").Value = "period_hours"
$ws.Range("D2").Value = 10
$ws.Columns.Item(4).ColumnWidth = 10.42

# --- 2. Insert "young_percentage" as new column I (after "initial_infected_percentage") ---
$ws.Columns.Item(9).Insert()
$ws.Range("I1").Value = "young_percentage"
$ws.Range("I2").Value = 0.8
# copy the formatting (style) from the neighbouring "initial_infected_percentage" cells
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- 3. Rename old "infection_probability" (now column J) to "infection_prob" ---
$ws.Range("J1").Value = "infection_prob"

# --- 4. Append two new trailing columns: K = reinfection_prob, L = vaccinated_infection_prob ---
$ws.Range("K1").Value = "reinfection_prob"
$ws.Range("K2").Value = 0.001

$ws.Range("L1").Value = "vaccinated_infection_prob"
$ws.Range("L2").Value = 0.003
$ws.Range("J1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("J2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- 5. Column widths for the two freshly appended columns ---
$ws.Columns.Item(11).ColumnWidth = 18.09
$ws.Columns.Item(12).ColumnWidth = 20.92

# --- 6. Sheet view: zoom + selection matching the new layout ---
$ws.Range("L5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 334
